$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 26, shifting the existing rows 26-32 down to 27-33.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly price-report entry.
$ws.Range("A26").Value2 = 7
$ws.Range("B26").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value2 = "Ñuble"
$ws.Range("D26").Value2 = 44637
$ws.Range("E26").Value2 = 16
$ws.Range("F26").Value2 = 100112040
$ws.Range("G26").Value2 = "Cilantro"
$ws.Range("H26").Value2 = "Sin especificar"
$ws.Range("I26").Value2 = "Primera"
$ws.Range("J26").Value2 = 120
$ws.Range("K26").Value2 = 550
$ws.Range("L26").Value2 = 600
$ws.Range("M26").Value2 = 575
$ws.Range("N26").Value2 = "`$/atado 0,5 a 1 kilo"
$ws.Range("O26").Value2 = "Provincia de Diguillín"
$ws.Range("P26").Value2 = 575
$ws.Range("Q26").Value2 = 1
$ws.Range("R26").Value2 = "Hortaliza"
